$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete rows 4-10 ("... tại HỆ THỐNG" entries) as a single block.
$ws.Range("A4:B10").EntireRow.Delete()

# After the above deletion, the old row 35 ("Tổng lương tại HỆ THỐNG")
# has shifted up to row 28. Delete it too.
$ws.Range("A28:B28").EntireRow.Delete()
